$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.283.35"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.45%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.058.42"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  +0.17%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "548.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "140.17"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.08%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.051.66"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("E9").Value = "  +0.56%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.44"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +2.13%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.82"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.555.86"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.15%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "63.245.49"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.48%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.060.61"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("E18").Value = "  -1.33%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.74"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.60%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "482.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.65%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.66"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("E22").Value = "  -0.61%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.23"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.51%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "80.59"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.58"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.04%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.75"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +4.00%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "26.04"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  -0.18%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +7.10%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.69"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.24%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "55.39"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.98"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "467.89"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.71%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0822"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.93%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0397"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.22%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.069.51"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.18%  "
$ws.Range("E41").Value = "  +0.13%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.26"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.05%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.39%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "28.14"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.255"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.08%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.06"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  -2.01%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0₃0509"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("E51").Value = "  +2.40%  "
